$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Divider row 342 ("THURSDAY") - copy styling/format from the existing divider row 313 ---
$ws.Range("A313:F313").Copy($ws.Range("A342:F342"))

# Row 343
$ws.Range("A343").Value = 'AV Shutdown'
$ws.Range("B343").Value = 42656
$ws.Range("C343").Value = '1700'
$ws.Range("D343").Value = 'FC'
$ws.Range("E343").Value = '305'
$ws.Range("F343").Value = 'LEAVE EQUIPMENT IN ROOM. Just log off PC and turn off projector. Lock room. Key for room in Founders 164 storeroom.'
$ws.Rows.Item(343).RowHeight = 30

# Row 344
$ws.Range("A344").Value = 'Pickup Mic'
$ws.Range("B344").Value = 42656
$ws.Range("C344").Value = '1600'
$ws.Range("D344").Value = 'VC'
$ws.Range("E344").Value = '001 -JCR'
$ws.Range("F344").Value = 'Pick up 3 desk mics, all cables and clips and NECK MIC ALSO. Return to Van 040 basement storeroom.'
$ws.Rows.Item(344).RowHeight = 30

# Row 345
$ws.Range("A345").Value = 'Pickup Small PA'
$ws.Range("B345").Value = 42656
$ws.Range("C345").Value = '1600'
$ws.Range("D345").Value = 'VC'
$ws.Range("E345").Value = '001-JCR'
$ws.Range("F345").Value = 'Pick up 2 small speakers and speaker cables, all matts and AC cords. Return to Van 040 storeroom.'
$ws.Rows.Item(345).RowHeight = 30

# Row 346
$ws.Range("A346").Value = 'Pickup PC'
$ws.Range("B346").Value = 42656
$ws.Range("C346").Value = '1600'
$ws.Range("D346").Value = 'VC'
$ws.Range("E346").Value = '001-JCR'
$ws.Range("F346").Value = 'Pick up roll in PC and Projector carts. Return to Vanier 040 basement storeroom. Key is in Founders 164 storeroom.'
$ws.Rows.Item(346).RowHeight = 30

# Row 347
$ws.Range("A347").Value = 'AV Shutdown'
$ws.Range("B347").Value = 42656
$ws.Range("C347").Value = '1800'
$ws.Range("D347").Value = 'R'
$ws.Range("E347").Value = 'N940'
$ws.Range("F347").Value = 'Senate Chamber - keys for room in Ross S120 storeroom. Turn off projector with remote on PC cart. Turn off PC. Turn off amplifier in back booth. Lock room.'
$ws.Rows.Item(347).RowHeight = 45

# Row 348
$ws.Range("A348").Value = 'AV Shutdown'
$ws.Range("B348").Value = 42656
$ws.Range("C348").Value = '1730'
$ws.Range("D348").Value = 'R'
$ws.Range("E348").Value = 'N102'
$ws.Range("F348").Value = 'Nat Taylor Cinema. Lock cinema all doors after shutdown.'

# Row 349
$ws.Range("A349").Value = 'Setup Skype Kit'
$ws.Range("B349").Value = 42656
$ws.Range("C349").Value = '1730'
$ws.Range("D349").Value = 'WC'
$ws.Range("E349").Value = '117'
$ws.Range("F349").Value = 'Set up Skype camera with tripod with built in PC in room. Skype camera is in Founders 164 storeroom. Tell prof. to stay with equipment until picked up. Tell Masi when end time is.'
$ws.Rows.Item(349).RowHeight = 45

# Row 350
$ws.Range("A350").Value = 'Pickup Skype Kit'
$ws.Range("B350").Value = 42656
$ws.Range("C350").Value = '2030'
$ws.Range("D350").Value = 'WC'
$ws.Range("E350").Value = '117'
$ws.Range("F350").Value = 'Pick up Skype camera and tripod. Log off PC and crestron and return skype camera and tripod to Founders 164 storeroom.'
$ws.Rows.Item(350).RowHeight = 30

# Row 351
$ws.Range("A351").Value = 'Other'
$ws.Range("B351").Value = 42656
$ws.Range("C351").Value = '1800'
$ws.Range("D351").Value = 'R'
$ws.Range("E351").Value = 'N102'
$ws.Range("F351").Value = 'Open doors to Nat Taylor Cinema. Allen key in S120 Ross storeroom.'
$ws.Rows.Item(351).RowHeight = 30

# Row 352
$ws.Range("A352").Value = 'AV Shutdown'
$ws.Range("B352").Value = 42656
$ws.Range("C352").Value = '2200'
$ws.Range("D352").Value = 'R'
$ws.Range("E352").Value = 'N102'
$ws.Range("F352").Value = 'Nat Taylor Cinema. Lock cinema all doors after shutdown.'

# Row 353
$ws.Range("A353").Value = 'Setup Mic'
$ws.Range("B353").Value = 42656
$ws.Range("C353").Value = '1630'
$ws.Range("D353").Value = 'VH'
$ws.Range("E353").Value = 'A'
$ws.Range("F353").Value = 'Set up 2 desk mics with mixer and neck microphone. Milk carton with mic cables, mics, mixer and desk stands is in Vari 1019 storeroom. Volume control on crestron - press mic icon and then ramp up " Podium" volume control.  Neck mic is in podium drawer.'
$ws.Rows.Item(353).RowHeight = 75

# Row 354
$ws.Range("A354").Value = 'Other'
$ws.Range("B354").Value = 42656
$ws.Range("C354").Value = '1630'
$ws.Range("D354").Value = 'VH'
$ws.Range("E354").Value = 'A'
$ws.Range("F354").Value = 'All equipment in milk crate in Vari 1019. Night Tech to meet Suzanne in room.'
$ws.Rows.Item(354).RowHeight = 30

# Row 355
$ws.Range("A355").Value = 'Pickup Mic'
$ws.Range("B355").Value = 42656
$ws.Range("C355").Value = '2200'
$ws.Range("D355").Value = 'VH'
$ws.Range("E355").Value = 'A'
$ws.Range("F355").Value = 'Pick up 2 mics, mic cables, 2 desk stands and mixer. Return to Vari 1019 storeroom. Wireless neck mic goes back to podium drawer.'
$ws.Rows.Item(355).RowHeight = 30

# --- View state: freeze header row, scroll viewport, set active cell/selection ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F358").Select()

